# Remove exercise-link placeholders from the class schedule that aren't
# ready yet ("update exercises on class schedule to remove before ready").
# These are all in column D ("In-Class Exercise") of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToClear = @(11, 12, 13, 14, 15, 18, 19, 20, 21, 22, 23, 24, 25, 29, 30)

foreach ($r in $rowsToClear) {
    $ws.Range("D$r").ClearContents()
}

$ws.Range("C10").Select()
